$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H15").Value = 241709.36
$ws.Range("I15").Value = 241709.36
$ws.Range("K15").Value = 725128.08
$ws.Range("M15").Value = -724959.08
$ws.Range("H38").Value = 300.25
$ws.Range("I38").Value = 62.833332
$ws.Range("J38").Value = 537.6667
$ws.Range("K38").Value = 188.499996
$ws.Range("L38").Value = 1613.0001
$ws.Range("M38").Value = 183.500004
$ws.Range("N38").Value = -2357.0001
$ws.Range("H76").Value = 3271197.2
$ws.Range("J76").Value = 3728.5715
$ws.Range("L76").Value = 3728.5715
$ws.Range("N76").Value = -4358.5715
$ws.Range("H79").Value = 3271197.2
$ws.Range("J79").Value = 3728.5715
$ws.Range("L79").Value = 3728.5715
$ws.Range("N79").Value = -5912.5715
$ws.Range("H125").Value = 37372372
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 37372372
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("M125").Value = 336351348
$ws.Range("N125").Value = -336356268
$ws.Range("H132").Value = 29854.086
$ws.Range("I132").Value = 30727.03
$ws.Range("J132").Value = 174
$ws.Range("K132").Value = 92181.09
$ws.Range("L132").Value = 522
$ws.Range("M132").Value = -89651.09
$ws.Range("N132").Value = -5582
$ws.Range("H135").Value = 2493.1924
$ws.Range("I135").Value = 2105.35
$ws.Range("J135").Value = 3786
$ws.Range("K135").Value = 18948.15
$ws.Range("L135").Value = 34074
$ws.Range("M135").Value = -16413.15
$ws.Range("N135").Value = -39144
# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 2080.3928
$ws.Range("I61").Value = 1707.1538
$ws.Range("J61").Value = 6932.5
$ws.Range("K61").Value = 1707.1538
$ws.Range("L61").Value = 6932.5
$ws.Range("M61").Value = -1495.1538
$ws.Range("N61").Value = -7356.5
$ws.Range("H80").Value = 99000
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 99000
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 6566.875
$ws.Range("I97").Value = 6964.6665
$ws.Range("J97").Value = 600
$ws.Range("K97").Value = 6964.6665
$ws.Range("L97").Value = 600
$ws.Range("M97").Value = -6468.6665
$ws.Range("N97").Value = -1592
$ws.Range("H101").Value = 19078
$ws.Range("J101").Value = 19078
$ws.Range("L101").Value = 19078
$ws.Range("N101").Value = -25568
$ws.Range("H132").Value = 3199.8438
$ws.Range("I132").Value = 2720
$ws.Range("J132").Value = 5279.1665
$ws.Range("K132").Value = 8160
$ws.Range("L132").Value = 15837.4995
$ws.Range("M132").Value = -5630
$ws.Range("N132").Value = -20897.4995
$ws.Range("H136").Value = 2080.3928
$ws.Range("I136").Value = 1707.1538
$ws.Range("J136").Value = 6932.5
$ws.Range("K136").Value = 5121.4614
$ws.Range("L136").Value = 20797.5
$ws.Range("M136").Value = -2571.4614
$ws.Range("N136").Value = -25897.5
# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 1143.1666
$ws.Range("I94").Value = 1106.9524
$ws.Range("J94").Value = 1396.6666
$ws.Range("K94").Value = 1106.9524
$ws.Range("L94").Value = 1396.6666
$ws.Range("M94").Value = -655.9523999999999
$ws.Range("N94").Value = -2298.6666
$ws.Range("H134").Value = 3905.2188
$ws.Range("I134").Value = 1997.2858
$ws.Range("J134").Value = 7547.636
$ws.Range("K134").Value = 5991.857400000001
$ws.Range("L134").Value = 22642.908
$ws.Range("M134").Value = -3456.857400000001
$ws.Range("N134").Value = -27712.908
# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H75").Value = 25000
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26996
$ws.Range("H78").Value = 25000
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84984
$ws.Range("H86").Value = 55558840
$ws.Range("J86").Value = 3939.2856
$ws.Range("L86").Value = 3939.2856
$ws.Range("N86").Value = -6185.2856
$ws.Range("H89").Value = 55558840
$ws.Range("J89").Value = 3939.2856
$ws.Range("L89").Value = 19696.428
$ws.Range("N89").Value = -30928.428
$ws.Range("H96").Value = 13483
$ws.Range("J96").Value = 13483
$ws.Range("L96").Value = 13483
$ws.Range("N96").Value = -18975
$ws.Range("H106").Value = 31866.666
$ws.Range("J106").Value = 31866.666
$ws.Range("L106").Value = 31866.666
$ws.Range("N106").Value = -34390.666
$ws.Range("H132").Value = 2882.2307
$ws.Range("I132").Value = 1513.1428
$ws.Range("K132").Value = 4539.428400000001
$ws.Range("M132").Value = -2009.428400000001
$ws.Range("H134").Value = 3295.2
$ws.Range("I134").Value = 1447.5333
$ws.Range("J134").Value = 6066.7
$ws.Range("K134").Value = 4342.5999
$ws.Range("L134").Value = 18200.1
$ws.Range("M134").Value = -1807.5999
$ws.Range("N134").Value = -23270.1
# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H31").Value = 1226
$ws.Range("I31").Value = 501
$ws.Range("J31").Value = 1951
$ws.Range("K31").Value = 1503
$ws.Range("L31").Value = 5853
$ws.Range("M31").Value = -1215
$ws.Range("N31").Value = -6429
$ws.Range("H129").Value = 1429.7222
$ws.Range("I129").Value = 1667
$ws.Range("J129").Value = 1278.7273
$ws.Range("K129").Value = 5001
$ws.Range("L129").Value = 3836.1819
$ws.Range("M129").Value = -1
$ws.Range("N129").Value = -13836.1819
# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H92").Value = 7225
$ws.Range("J92").Value = 7225
$ws.Range("L92").Value = 7225
$ws.Range("N92").Value = -10969
$ws.Range("H132").Value = 3117.25
$ws.Range("I132").Value = 2843.9412
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 8531.8236
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -6001.8236
$ws.Range("N132").Value = -19058
# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H68").Value = 2069
$ws.Range("I68").Value = 1927.1428
$ws.Range("J68").Value = 2400
$ws.Range("K68").Value = 1927.1428
$ws.Range("L68").Value = 2400
$ws.Range("M68").Value = -1178.1428
$ws.Range("N68").Value = -3898
$ws.Range("H70").Value = 19800
$ws.Range("J70").Value = 19800
$ws.Range("L70").Value = 19800
$ws.Range("N70").Value = -20340
$ws.Range("H71").Value = 2069
$ws.Range("I71").Value = 1927.1428
$ws.Range("J71").Value = 2400
$ws.Range("K71").Value = 9635.714
$ws.Range("L71").Value = 12000
$ws.Range("M71").Value = -5891.714
$ws.Range("N71").Value = -19488
$ws.Range("H73").Value = 19800
$ws.Range("J73").Value = 19800
$ws.Range("L73").Value = 19800
$ws.Range("N73").Value = -21672
$ws.Range("H82").Value = 1412.2667
$ws.Range("I82").Value = 1152
$ws.Range("J82").Value = 2453.3333
$ws.Range("K82").Value = 1152
$ws.Range("L82").Value = 2453.3333
$ws.Range("M82").Value = -791
$ws.Range("N82").Value = -3175.3333
$ws.Range("H85").Value = 1412.2667
$ws.Range("I85").Value = 1152
$ws.Range("J85").Value = 2453.3333
$ws.Range("K85").Value = 1152
$ws.Range("L85").Value = 2453.3333
$ws.Range("M85").Value = 96
$ws.Range("N85").Value = -4949.3333
$ws.Range("H98").Value = 29355
$ws.Range("J98").Value = 29355
$ws.Range("L98").Value = 29355
$ws.Range("N98").Value = -35345
$ws.Range("H106").Value = 16513.75
$ws.Range("J106").Value = 16513.75
$ws.Range("L106").Value = 16513.75
$ws.Range("N106").Value = -19037.75
$ws.Range("H132").Value = 3312.3684
$ws.Range("I132").Value = 1995.28
$ws.Range("J132").Value = 5845.231
$ws.Range("K132").Value = 5985.84
$ws.Range("L132").Value = 17535.693
$ws.Range("M132").Value = -3455.84
$ws.Range("N132").Value = -22595.693
# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H97").Value = 29293.334
$ws.Range("J97").Value = 29293.334
$ws.Range("L97").Value = 29293.334
$ws.Range("N97").Value = -31275.334
$ws.Range("H103").Value = 273850.5
$ws.Range("J103").Value = 273850.5
$ws.Range("L103").Value = 273850.5
$ws.Range("N103").Value = -276194.5
$ws.Range("H105").Value = 19307.5
$ws.Range("J105").Value = 19307.5
$ws.Range("L105").Value = 19307.5
$ws.Range("N105").Value = -26295.5
$ws.Range("H122").Value = 1339.8889
$ws.Range("I122").Value = 1151
$ws.Range("J122").Value = 1491
$ws.Range("K122").Value = 3453
$ws.Range("L122").Value = 4473
$ws.Range("M122").Value = -1003
$ws.Range("N122").Value = -9373
$ws.Range("H132").Value = 2913.9783
$ws.Range("I132").Value = 2935
$ws.Range("J132").Value = 2838.3
$ws.Range("K132").Value = 8805
$ws.Range("L132").Value = 8514.900000000001
$ws.Range("M132").Value = -6275
$ws.Range("N132").Value = -13574.9
$ws.Range("H136").Value = 1800.5217
$ws.Range("I136").Value = 1077.1765
$ws.Range("J136").Value = 3850
$ws.Range("K136").Value = 3231.5295
$ws.Range("L136").Value = 11550
$ws.Range("M136").Value = -681.5295000000001
$ws.Range("N136").Value = -16650
